$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DATE column (B) for the log rows (2-43) moves from 2024-11-02 to
# 2024-11-09 (one week later), matching the aws cloudwatch re-run noted in
# the commit message. The dates are stored as literal text in the sheet
# (not real date values), so write them back as text too: prefixing with
# an apostrophe stops Excel from reinterpreting "2024-11-09" as a date
# serial number. That apostrophe trick marks the cell with a "quote
# prefix" format, so restore the original cell formatting afterwards by
# copying it over from an adjacent, untouched column in the same rows.

$ws.Range("C2:C43").Copy() | Out-Null
$ws.Range("B2:B43").Value = "'2024-11-09"

$ws.Range("C2:C43").Copy() | Out-Null
$ws.Range("B2:B43").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
